# aggiornamento fino a 27/05
# Appends 14 new daily rows (rows 256-269) to the single data sheet, covering
# dates 2021-05-14 through 2021-05-27 (Excel serials 44330-44343).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data: date-serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti
$data = @(
  @(44330, 0, 2, 47.65308553728854),
  @(44331, 1, 1, 23.82654276864427),
  @(44332, 0, 1, 23.82654276864427),
  @(44333, 0, 1, 23.82654276864427),
  @(44334, 0, 1, 23.82654276864427),
  @(44335, 0, 1, 23.82654276864427),
  @(44336, 0, 1, 23.82654276864427),
  @(44337, 0, 1, 23.82654276864427),
  @(44338, 0, 0, 0),
  @(44339, 0, 0, 0),
  @(44340, 0, 0, 0),
  @(44341, 0, 0, 0),
  @(44342, 0, 0, 0),
  @(44343, 0, 0, 0)
)

$startRow = 256
$endRow = $startRow + $data.Count - 1

# Column A keeps the same date style as the rest of the column (centered,
# bordered, custom date/time number format) -- copy that formatting down
# from the last existing row rather than re-creating a style definition.
$ws.Range("A255").Copy()
$ws.Range("A" + $startRow + ":A" + $endRow).PasteSpecial(-4122)

$r = $startRow
foreach ($row in $data) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  $r = $r + 1
}

Write-Output "Added rows $startRow to $endRow"
